$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.045.15"
$ws.Range("E2").Value = "  -3.70%  "

# Row 3
$ws.Range("D3").Value = "1.639.74"
$ws.Range("E3").Value = "  -3.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3876"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.69%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3837"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.00%  "

# Row 9
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.95%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.341"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08419"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.06%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.078"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001273"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.75%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.434"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.82%  "

# Row 17
$ws.Range("D17").Value = "1.644.91"
$ws.Range("E17").Value = "  -5.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06934"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.893"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.56"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "24.039.78"
$ws.Range("E24").Value = "  -3.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.667"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.55%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.29%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.32%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.718"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.91%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "141.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.245"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.454"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.06%  "

# Row 33
$ws.Range("D33").Value = "1.822.83"
$ws.Range("E33").Value = "  -7.48%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.104"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07967"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02902"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9549"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.75%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09189"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.20%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.460"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.907"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7549"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.54%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6847"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.464"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.079"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08335"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.70%  "
